# Update the "Förändrad" (Modified) date column (C) for rows 2-37
# from 45658 (2025-01-01) to 45659 (2025-01-02).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 37; $row++) {
    $ws.Cells.Item($row, 3).Value = 45659
}
